# Add "Udon Don Bar" as a new store row (with its opening hours) into
# Sheet1, directly above the existing "Waa Cow!" row, keeping the sheet
# alphabetically sorted by Store name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 26. This shifts "Waa Cow!" (and the
# trailing blank spacer row) down by one row, and the new row inherits
# the cell styles of the row above it (row 25).
$ws.Rows.Item(26).Insert()

# Populate the new row with the Udon Don Bar data.
$ws.Cells.Item(26, 1).Value = "Udon Don Bar"
$ws.Cells.Item(26, 2).Value = "Food & Beverages"
$ws.Cells.Item(26, 3).Value = "Japanese Cuisine"
$ws.Cells.Item(26, 4).Value = "No"
$ws.Cells.Item(26, 5).Value = "Town Plaza"
$ws.Cells.Item(26, 6).Value = "1100-2200"
$ws.Cells.Item(26, 7).Value = "1100-2200"
$ws.Cells.Item(26, 8).Value = "1100-2200"
$ws.Cells.Item(26, 9).Value = "1100-2200"
$ws.Cells.Item(26, 10).Value = "1100-2200"
$ws.Cells.Item(26, 11).Value = "1100-2200"
$ws.Cells.Item(26, 12).Value = "Closed"
$ws.Cells.Item(26, 13).Value = "Closed"

# Match the row height used by the other data rows.
$ws.Rows.Item(26).RowHeight = 13.2

# The data range grew by one row (now ends at row 27 instead of 26):
# update the hidden defined names that track the filter database range.
$names = $wb.Names
$names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$U`$27"
$names.Item(2).RefersTo = "=Sheet1!`$A`$1:`$F`$27"

# Re-apply the AutoFilter so it covers the new, larger range.
$ws.AutoFilterMode = $False | Out-Null
$ws.Range("A1:U27").AutoFilter() | Out-Null

# Restore the view state (active cell / scroll position) as left by the
# editor.
$ws.Range("J17").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 8
